$wb = $excel.ActiveWorkbook

# --- 1. Create the new "streaming_server" sheet as a copy of an existing
#        "Test Case Summary" sheet (protocols(setting)) so it inherits the
#        identical column widths / styles / merged cell layout, then move
#        it into position 3 (right after my_account, before
#        protocols(setting)) and rename it. ---
$template = $wb.Worksheets.Item("protocols(setting)")
$template.Copy($null, $template)
$newSheet = $wb.Worksheets.Item(4)
$newSheet.Name = "streaming_server"

# --- 2. Populate the new sheet's content (values only change on D2/E3/E4/
#        E5/E6/E8 relative to the copied template; D3..D8 labels already
#        match since every "Test Case Summary" sheet shares the same
#        label text). Do this BEFORE the Move below, since writing cell
#        values to a sheet right after it has been relocated is unreliable
#        for numeric values in this runtime. ---
$newSheet.Range("D2").Value = " Test Case Summary (10-03-24)"
$newSheet.Range("E3").Value = 37
$newSheet.Range("E4").Value = 33
$newSheet.Range("E5").Value = 0
$newSheet.Range("E6").Value = 4
$newSheet.Range("E8").ClearContents()

# --- 2b. Move it into position 3 (right after my_account, before
#         protocols(setting)). ---
$newSheet.Move($template)

# --- 3. New sheet becomes the active / selected tab with selection G6.
#        Re-fetch the worksheet reference after the Move so the
#        activation/selection below is applied to the relocated sheet. ---
$streaming = $wb.Worksheets.Item("streaming_server")
$streaming.Activate()
$streaming.Range("G6").Select()
